$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132 (sheet ALC)
$ws.Range("H132").Value = 2371.375
$ws.Range("I132").Value = 2286.3872
$ws.Range("J132").Value = 5006
$ws.Range("K132").Value = 6859.1616
$ws.Range("L132").Value = 15018
$ws.Range("M132").Value = -4329.1616
$ws.Range("N132").Value = -20078

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (sheet ARM)
$ws.Range("H45").Value = 2083.1714
$ws.Range("I45").Value = 2218.6
$ws.Range("J45").Value = 1744.6
$ws.Range("K45").Value = 2218.6
$ws.Range("L45").Value = 1744.6
$ws.Range("M45").Value = -1841.6
$ws.Range("N45").Value = -2498.6

# Row 110 (sheet ARM)
$ws.Range("H110").Value = 24063.6
$ws.Range("I110").Value = 29925.666
$ws.Range("K110").Value = 29925.666
$ws.Range("M110").Value = -27880.666

# Row 112 (sheet ARM)
$ws.Range("H112").Value = 20157.4
$ws.Range("J112").Value = 20157.4
$ws.Range("L112").Value = 20157.4
$ws.Range("N112").Value = -23111.4

# Row 132 (sheet ARM)
$ws.Range("H132").Value = 33419030
$ws.Range("J132").Value = 125301530
$ws.Range("L132").Value = 375904590
$ws.Range("N132").Value = -375909650

$ws = $wb.Worksheets.Item("BSM")
# Row 134 (sheet BSM)
$ws.Range("H134").Value = 2502432
$ws.Range("I134").Value = 2779660.5
$ws.Range("K134").Value = 8338981.5
$ws.Range("M134").Value = -8336446.5

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (sheet CRP)
$ws.Range("H7").Value = 298.83334
$ws.Range("I7").Value = 101.57143
$ws.Range("J7").Value = 424.36365
$ws.Range("K7").Value = 101.57143
$ws.Range("L7").Value = 424.36365
$ws.Range("M7").Value = 11.42856999999999
$ws.Range("N7").Value = -650.36365

# Row 99 (sheet CRP)
$ws.Range("H99").Value = 8958.433000000001
$ws.Range("J99").Value = 10461.538
$ws.Range("L99").Value = 10461.538
$ws.Range("N99").Value = -13457.538

# Row 126 (sheet CRP)
$ws.Range("H126").Value = 8958.433000000001
$ws.Range("J126").Value = 10461.538
$ws.Range("L126").Value = 31384.614
$ws.Range("N126").Value = -36324.614

# Row 134 (sheet CRP)
$ws.Range("H134").Value = 2459.5334
$ws.Range("I134").Value = 2326.8333
$ws.Range("K134").Value = 6980.499899999999
$ws.Range("M134").Value = -4445.499899999999

$ws = $wb.Worksheets.Item("CUL")
# Row 9 (sheet CUL)
$ws.Range("H9").Value = 762.5
$ws.Range("I9").Value = 516.6667
$ws.Range("K9").Value = 1550.0001
$ws.Range("M9").Value = -1326.0001

# Row 129 (sheet CUL)
$ws.Range("H129").Value = 2798.158
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 2798.158
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 8394.474
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -18394.474

# Row 131 (sheet CUL)
$ws.Range("H131").Value = 1429.3684
$ws.Range("J131").Value = 1694.5
$ws.Range("L131").Value = 5083.5
$ws.Range("N131").Value = -15163.5

# Row 132 (sheet CUL)
$ws.Range("H132").Value = 2088554.9
$ws.Range("I132").Value = 1754.091
$ws.Range("J132").Value = 3181641
$ws.Range("K132").Value = 15786.819
$ws.Range("L132").Value = 28634769
$ws.Range("M132").Value = -13256.819
$ws.Range("N132").Value = -28639829

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (sheet LTW)
$ws.Range("H7").Value = 4717.9
$ws.Range("I7").Value = 5141.273
$ws.Range("K7").Value = 5141.273
$ws.Range("M7").Value = -5029.273

# Row 22 (sheet LTW)
$ws.Range("H22").Value = 2861.3333
$ws.Range("I22").Value = 1808.6666
$ws.Range("K22").Value = 1808.6666
$ws.Range("M22").Value = -1513.6666

# Row 27 (sheet LTW)
$ws.Range("H27").Value = 2861.3333
$ws.Range("I27").Value = 1808.6666
$ws.Range("K27").Value = 1808.6666
$ws.Range("M27").Value = -1701.6666

# Row 40 (sheet LTW)
$ws.Range("H40").Value = 3146.4546
$ws.Range("I40").Value = 3041.1
$ws.Range("J40").Value = 4200
$ws.Range("K40").Value = 3041.1
$ws.Range("L40").Value = 4200
$ws.Range("M40").Value = -2905.1
$ws.Range("N40").Value = -4472

# Row 61 (sheet LTW)
$ws.Range("H61").Value = 5188.7036
$ws.Range("I61").Value = 3528.5
$ws.Range("K61").Value = 3528.5
$ws.Range("M61").Value = -3326.5

# Row 63 (sheet LTW)
$ws.Range("H63").Value = 78000
$ws.Range("I63").Value = 78000
$ws.Range("K63").Value = 78000
$ws.Range("M63").Value = -77251

# Row 66 (sheet LTW)
$ws.Range("H66").Value = 78000
$ws.Range("I66").Value = 78000
$ws.Range("K66").Value = 234000
$ws.Range("M66").Value = -230256

# Row 68 (sheet LTW)
$ws.Range("H68").Value = 2546.3333
$ws.Range("I68").Value = 1953.1666
$ws.Range("K68").Value = 1953.1666
$ws.Range("M68").Value = -1204.1666

# Row 71 (sheet LTW)
$ws.Range("H71").Value = 2546.3333
$ws.Range("I71").Value = 1953.1666
$ws.Range("K71").Value = 9765.833000000001
$ws.Range("M71").Value = -6021.833000000001

# Row 93 (sheet LTW)
$ws.Range("H93").Value = 2483.111
$ws.Range("I93").Value = 2079.0833
$ws.Range("K93").Value = 2079.0833
$ws.Range("M93").Value = -831.0832999999998

# Row 113 (sheet LTW)
$ws.Range("H113").Value = 5188.7036
$ws.Range("I113").Value = 3528.5
$ws.Range("K113").Value = 3528.5
$ws.Range("M113").Value = -1358.5

# Row 122 (sheet LTW)
$ws.Range("H122").Value = 3753
$ws.Range("I122").Value = 2998.75
$ws.Range("K122").Value = 8996.25
$ws.Range("M122").Value = -6546.25

# Row 126 (sheet LTW)
$ws.Range("H126").Value = 4717.9
$ws.Range("I126").Value = 5141.273
$ws.Range("K126").Value = 15423.819
$ws.Range("M126").Value = -12953.819

# Row 136 (sheet LTW)
$ws.Range("H136").Value = 2330.9285
$ws.Range("I136").Value = 2371
$ws.Range("J136").Value = 688
$ws.Range("K136").Value = 7113
$ws.Range("L136").Value = 2064
$ws.Range("M136").Value = -4563
$ws.Range("N136").Value = -7164

$ws = $wb.Worksheets.Item("WVR")
# Row 2 (sheet WVR)
$ws.Range("H2").Value = 202500
$ws.Range("I2").Value = 202500
$ws.Range("K2").Value = 202500
$ws.Range("M2").Value = -202388

# Row 4 (sheet WVR)
$ws.Range("H4").Value = 13456.5
$ws.Range("I4").Value = 18217
$ws.Range("K4").Value = 18217
$ws.Range("M4").Value = -18104

# Row 6 (sheet WVR)
$ws.Range("H6").Value = 4000
$ws.Range("I6").Value = 4000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 4000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -3885
$ws.Range("N6").ClearContents()

# Row 96 (sheet WVR)
$ws.Range("H96").Value = 4600
$ws.Range("J96").Value = 9000
$ws.Range("L96").Value = 9000
$ws.Range("N96").Value = -11746

# Row 100 (sheet WVR)
$ws.Range("H100").Value = 48097532
$ws.Range("I100").Value = 59413980
$ws.Range("K100").Value = 118827960
$ws.Range("M100").Value = -118827419

# Row 122 (sheet WVR)
$ws.Range("H122").Value = 55612516
$ws.Range("I122").Value = 66734332
$ws.Range("J122").Value = 3433.3333
$ws.Range("K122").Value = 200202996
$ws.Range("L122").Value = 10299.9999
$ws.Range("M122").Value = -200200546
$ws.Range("N122").Value = -15199.9999

# Row 126 (sheet WVR)
$ws.Range("H126").Value = 11072
$ws.Range("I126").Value = 16899.8
$ws.Range("J126").Value = 3787.25
$ws.Range("K126").Value = 50699.39999999999
$ws.Range("L126").Value = 11361.75
$ws.Range("M126").Value = -48229.39999999999
$ws.Range("N126").Value = -16301.75

# Row 132 (sheet WVR)
$ws.Range("H132").Value = 1786.3846
$ws.Range("I132").Value = 1574.4166
$ws.Range("K132").Value = 4723.2498
$ws.Range("M132").Value = -2193.2498
